$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price/Volume columns so that numeric-looking
# strings (e.g. "553.08", "67.895.37") are kept as literal text, matching the
# workbooks original inline-string storage instead of being auto-coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "67.895.37"
$ws.Range("E2").Value = "  -1.28%  "
$ws.Range("D3").Value = "2.398.90"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "553.08"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "158.10"
$ws.Range("E6").Value = "  -2.52%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.157"
$ws.Range("E9").Value = "  +3.78%  "
$ws.Range("E10").Value = "  -1.40%  "
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("D12").Value = "4.71"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "67.770.53"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("E15").Value = "  -3.37%  "
$ws.Range("D16").Value = "10.27"
$ws.Range("E16").Value = "  -4.42%  "
$ws.Range("D17").Value = "329.44"
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "6.79"
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "65.51"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").Value = "3.62"
$ws.Range("E23").Value = "  -1.93%  "
$ws.Range("D24").Value = "8.03"
$ws.Range("E24").Value = "  -1.61%  "
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").Value = "7.03"
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").Value = "418.77"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D29").Value = "1.13"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").Value = "157.16"
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "18.97"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "17.59"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").Value = "0.292"
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("D37").Value = "4.22"
$ws.Range("E37").Value = "  -5.05%  "
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("E39").Value = "  -4.49%  "
$ws.Range("D40").Value = "3.27"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "128.45"
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("E42").Value = "  -7.05%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "0.474"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "0.552"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "0.0913"
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -6.07%  "
$ws.Range("D49").Value = "16.36"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E50").Value = "  -6.53%  "
$ws.Range("D51").Value = "0.0426"
$ws.Range("E51").Value = "  -0.35%  "
